# MasterExecutor_Sanity.xlsx - "Temp changed MAster KFP"
# Flip the Runmode column (E) from "Yes" to "No" for every test case row
# that was still set to "Yes" (rows already "No" are left untouched),
# and update the sheet's current selection to the Runmode column.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("MasterExecutor")

# Rows (data rows 3-29, excluding the header row 1 and row 2 / the rows
# that were already "No") whose Runmode value changes from Yes -> No.
$rowsToFlip = @(3,4,5,6,7,8,9,10,11,12,13,14,15,17,18,21,22,23,25,28,29)

foreach ($r in $rowsToFlip) {
    $ws.Cells.Item($r, 5).Value = "No"
}

# Match the workbook's recorded selection after the edit: the active
# sheet, selecting the Runmode column from row 3 down to the last row.
$ws.Activate()
$ws.Range("E3:E29").Select()
